$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregar cuenta 70106 al grupo RECURRENTE_<_100K (celda B5)
$current = $ws.Range("B5").Value()
$ws.Range("B5").Value = $current + ".70106"

# Actualizar la vista de la hoja: desplazar a B2 y seleccionar B6
$excel.Goto($ws.Range("B2"), $true)
$null = $ws.Range("B6").Select()
